$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "69.744.38"
Set-TextValue "E2" "  -0.12%  "

Set-TextValue "D3" "3.682.38"
Set-TextValue "E3" "  -0.61%  "

Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.11%  "

Set-TextValue "D5" "651.95"
Set-TextValue "E5" "  -3.89%  "

Set-TextValue "D6" "161.22"
Set-TextValue "E6" "  -0.89%  "

Set-TextValue "E7" "  +0.03%  "

Set-TextValue "D8" "0.498"
Set-TextValue "E8" "  +0.07%  "

Set-TextValue "E9" "  -2.30%  "

Set-TextValue "D10" "7.16"
Set-TextValue "E10" "  +0.38%  "

Set-TextValue "D11" "0.442"
Set-TextValue "E11" "  -0.50%  "

Set-TextValue "E12" "  -2.14%  "

Set-TextValue "D13" "4.303.09"
Set-TextValue "E13" "  -0.67%  "

Set-TextValue "D14" "32.70"
Set-TextValue "E14" "  -0.80%  "

Set-TextValue "D15" "3.670.16"
Set-TextValue "E15" "  -1.25%  "

Set-TextValue "D16" "69.741.52"
Set-TextValue "E16" "  -0.07%  "

Set-TextValue "E17" "  +0.66%  "

Set-TextValue "D18" "6.53"
Set-TextValue "E18" "  +0.26%  "

Set-TextValue "D19" "15.93"
Set-TextValue "E19" "  -1.28%  "

Set-TextValue "D20" "10.37"
Set-TextValue "E20" "  +5.33%  "

Set-TextValue "D21" "470.47"
Set-TextValue "E21" "  -0.54%  "

Set-TextValue "D22" "0.655"
Set-TextValue "E22" "  +0.09%  "

Set-TextValue "D23" "79.83"
Set-TextValue "E23" "  -0.93%  "

Set-TextValue "D24" "3.826.92"
Set-TextValue "E24" "  -0.70%  "

Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.00%  "

Set-TextValue "E26" "  -1.75%  "

Set-TextValue "D27" "11.15"
Set-TextValue "E27" "  +0.98%  "

Set-TextValue "E28" "  -4.35%  "

Set-TextValue "D29" "2.65"
Set-TextValue "E29" "  -2.29%  "

Set-TextValue "E30" "  -3.49%  "

Set-TextValue "B31" "ImmutableX"
Set-TextValue "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "1.99"
Set-TextValue "E31" "  -2.05%  "

Set-TextValue "B32" "Binance-PegBSC-USD"
Set-TextValue "C32" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  -0.04%  "

Set-TextValue "D33" "0.166"
Set-TextValue "E33" "  +0.59%  "

Set-TextValue "D34" "26.74"
Set-TextValue "E34" "  -0.76%  "

Set-TextValue "D35" "6.44"
Set-TextValue "E35" "  -3.01%  "

Set-TextValue "D36" "3.675.95"
Set-TextValue "E36" "  -0.52%  "

Set-TextValue "D37" "8.38"
Set-TextValue "E37" "  -2.53%  "

Set-TextValue "D39" "5.88"
Set-TextValue "E39" "  -5.28%  "

Set-TextValue "D40" "178.22"
Set-TextValue "E40" "  +5.53%  "

Set-TextValue "E41" "  -0.13%  "

Set-TextValue "D42" "0.0895"
Set-TextValue "E42" "  -1.52%  "

Set-TextValue "D43" "2.18"
Set-TextValue "E43" "  -2.37%  "

Set-TextValue "E44" "  -1.73%  "

Set-TextValue "D45" "46.82"
Set-TextValue "E45" "  -0.39%  "

Set-TextValue "D46" "29.09"
Set-TextValue "E46" "  +3.33%  "

Set-TextValue "E47" "  -0.90%  "

Set-TextValue "E48" "  -5.02%  "

Set-TextValue "D49" "7.86"
Set-TextValue "E49" "  -0.95%  "

Set-TextValue "E50" "  -5.12%  "

Set-TextValue "D51" "1.05"
Set-TextValue "E51" "  -6.13%  "
